$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J19").Value = 1975.1538
$ws.Range("N19").Value = -2325.1538
$ws.Range("L19").Value = 1975.1538
$ws.Range("H19").Value = 1658.8
$ws.Range("I62").Value = 2952.875
$ws.Range("K62").Value = 2952.875
$ws.Range("H62").Value = 3176.4167
$ws.Range("M62").Value = -2328.875
$ws.Range("M65").Value = -11644.375
$ws.Range("K65").Value = 14764.375
$ws.Range("H65").Value = 3176.4167
$ws.Range("I65").Value = 2952.875
$ws.Range("H125").Value = 149717.14
$ws.Range("I125").Value = 207704
$ws.Range("K125").Value = 1869336
$ws.Range("M125").Value = -1866876
$ws.Range("J138").Value = 2480.9167
$ws.Range("L138").Value = 7442.750100000001
$ws.Range("N138").Value = -17722.7501
$ws.Range("M138").Value = -2000
$ws.Range("K138").Value = 7140
$ws.Range("I138").Value = 2380
$ws.Range("H138").Value = 2460.7334

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4236.3564
$ws.Range("I32").Value = 1943.1805
$ws.Range("N32").Value = -15817.6
$ws.Range("M32").Value = -1656.1805
$ws.Range("L32").Value = 15243.6
$ws.Range("K32").Value = 1943.1805
$ws.Range("J32").Value = 15243.6
$ws.Range("J45").Value = 31501250
$ws.Range("H45").Value = 14002269
$ws.Range("N45").Value = -31502004
$ws.Range("L45").Value = 31501250
$ws.Range("J74").Value = 4160.5
$ws.Range("K74").Value = 72898.57000000001
$ws.Range("M74").Value = -72024.57000000001
$ws.Range("H74").Value = 52277.15
$ws.Range("L74").Value = 4160.5
$ws.Range("N74").Value = -5908.5
$ws.Range("I74").Value = 72898.57000000001
$ws.Range("I77").Value = 72898.57000000001
$ws.Range("M77").Value = -360124.85
$ws.Range("K77").Value = 364492.85
$ws.Range("L77").Value = 20802.5
$ws.Range("J77").Value = 4160.5
$ws.Range("H77").Value = 52277.15
$ws.Range("N77").Value = -29538.5
$ws.Range("M97").Value = -513.1429000000001
$ws.Range("H97").Value = 896
$ws.Range("I97").Value = 1009.1429
$ws.Range("K97").Value = 1009.1429
$ws.Range("K110").Value = 886
$ws.Range("I110").Value = 886
$ws.Range("M110").Value = 1159
$ws.Range("H110").Value = 1152.1111
$ws.Range("K132").Value = 6541.2855
$ws.Range("H132").Value = 2341.2856
$ws.Range("M132").Value = -4011.2855
$ws.Range("N132").Value = -13049
$ws.Range("L132").Value = 7989
$ws.Range("I132").Value = 2180.4285
$ws.Range("J132").Value = 2663
$ws.Range("J138").Value = 66499.5
$ws.Range("L138").Value = 66499.5
$ws.Range("N138").Value = -76779.5
$ws.Range("H138").Value = 66499.5
$ws.Range("H140").Value = 74660.336
$ws.Range("L140").Value = 74660.336
$ws.Range("J140").Value = 74660.336
$ws.Range("N140").Value = -85020.336

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("L20").Value = 1425
$ws.Range("N20").Value = -1919
$ws.Range("H20").Value = 1313
$ws.Range("J20").Value = 1425
$ws.Range("K86").Value = 4357.3
$ws.Range("J86").Value = 5029
$ws.Range("I86").Value = 4357.3
$ws.Range("H86").Value = 4633.8823
$ws.Range("L86").Value = 5029
$ws.Range("N86").Value = -7275
$ws.Range("M86").Value = -3234.3
$ws.Range("I89").Value = 4357.3
$ws.Range("H89").Value = 4633.8823
$ws.Range("K89").Value = 21786.5
$ws.Range("J89").Value = 5029
$ws.Range("N89").Value = -36377
$ws.Range("M89").Value = -16170.5
$ws.Range("L89").Value = 25145
$ws.Range("H105").Value = 129781.125
$ws.Range("N105").Value = -9544
$ws.Range("J105").Value = 6050
$ws.Range("I105").Value = 335999.66
$ws.Range("K105").Value = 335999.66
$ws.Range("M105").Value = -334252.66
$ws.Range("L105").Value = 6050
$ws.Range("N107").Value = -9395
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 5555
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 5555
$ws.Range("H107").Value = 5555
$ws.Range("N110").Value = -89602
$ws.Range("L110").Value = 81422
$ws.Range("J110").Value = 81422
$ws.Range("H110").Value = 81422
$ws.Range("M107").ClearContents()

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M22").Value = -1506.5714
$ws.Range("I22").Value = 1856.5714
$ws.Range("H22").Value = 2302.4707
$ws.Range("K22").Value = 1856.5714
$ws.Range("I31").Value = 1945.7
$ws.Range("K31").Value = 1945.7
$ws.Range("H31").Value = 2816.1482
$ws.Range("M31").Value = -1650.7
$ws.Range("M34").Value = -1743.7
$ws.Range("H34").Value = 2816.1482
$ws.Range("I34").Value = 1945.7
$ws.Range("K34").Value = 1945.7
$ws.Range("J86").Value = 8330
$ws.Range("H86").Value = 3973586
$ws.Range("L86").Value = 8330
$ws.Range("N86").Value = -10576
$ws.Range("H89").Value = 3973586
$ws.Range("J89").Value = 8330
$ws.Range("N89").Value = -52882
$ws.Range("L89").Value = 41650

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1102.4546
$ws.Range("K5").Value = 3313.5
$ws.Range("I5").Value = 1104.5
$ws.Range("M5").Value = -3201.5
$ws.Range("K74").Value = 21060
$ws.Range("M74").Value = -19999
$ws.Range("H74").Value = 11759.75
$ws.Range("I74").Value = 7020
$ws.Range("I77").Value = 7020
$ws.Range("M77").Value = -57876
$ws.Range("K77").Value = 63180
$ws.Range("H77").Value = 11759.75
$ws.Range("M97").Value = 152.5
$ws.Range("H97").Value = 121.07143
$ws.Range("I97").Value = 114.5
$ws.Range("K97").Value = 343.5
$ws.Range("J114").Value = 11837.625
$ws.Range("N114").Value = -42020.875
$ws.Range("L114").Value = 35512.875
$ws.Range("H114").Value = 11837.625
$ws.Range("H121").Value = 1991.1538
$ws.Range("L121").Value = 12000
$ws.Range("J121").Value = 4000
$ws.Range("M121").Value = -896.875
$ws.Range("I121").Value = 735.625
$ws.Range("K121").Value = 2206.875
$ws.Range("N121").Value = -14620
$ws.Range("M122").Value = -5949.9997
$ws.Range("K122").Value = 8399.9997
$ws.Range("J122").Value = 1500
$ws.Range("H122").Value = 1160
$ws.Range("N122").Value = -18400
$ws.Range("L122").Value = 13500
$ws.Range("I122").Value = 933.3333
$ws.Range("K132").Value = 9549
$ws.Range("H132").Value = 3642.4167
$ws.Range("M132").Value = -7019
$ws.Range("N132").Value = -49458.125
$ws.Range("L132").Value = 44398.125
$ws.Range("I132").Value = 1061
$ws.Range("J132").Value = 4933.125
$ws.Range("M135").Value = -7405.5
$ws.Range("I135").Value = 1104.5
$ws.Range("K135").Value = 9940.5
$ws.Range("H135").Value = 1102.4546

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L52").Value = 14028.667
$ws.Range("H52").Value = 14171.375
$ws.Range("N52").Value = -14546.667
$ws.Range("J52").Value = 14028.667
$ws.Range("M122").Value = -3723.4375
$ws.Range("K122").Value = 6173.4375
$ws.Range("J122").Value = 2873.8
$ws.Range("H122").Value = 2252.0952
$ws.Range("N122").Value = -13521.4
$ws.Range("L122").Value = 8621.400000000001
$ws.Range("I122").Value = 2057.8125
$ws.Range("K132").Value = 10209.8334
$ws.Range("H132").Value = 4229.1924
$ws.Range("M132").Value = -7679.8334
$ws.Range("N132").Value = -23322.5
$ws.Range("L132").Value = 18262.5
$ws.Range("I132").Value = 3403.2778
$ws.Range("J132").Value = 6087.5

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M7").Value = -2343.2
$ws.Range("H7").Value = 3445.95
$ws.Range("K7").Value = 2455.2
$ws.Range("I7").Value = 2455.2
$ws.Range("N16").Value = -4431.75
$ws.Range("M16").Value = -2370.077
$ws.Range("J16").Value = 4091.75
$ws.Range("K16").Value = 2540.077
$ws.Range("L16").Value = 4091.75
$ws.Range("H16").Value = 2905.1765
$ws.Range("I16").Value = 2540.077
$ws.Range("L22").Value = 757.4167
$ws.Range("M22").Value = -605
$ws.Range("N22").Value = -1347.4167
$ws.Range("J22").Value = 757.4167
$ws.Range("I22").Value = 900
$ws.Range("H22").Value = 822.2273
$ws.Range("K22").Value = 900
$ws.Range("L27").Value = 757.4167
$ws.Range("K27").Value = 900
$ws.Range("I27").Value = 900
$ws.Range("N27").Value = -971.4167
$ws.Range("H27").Value = 822.2273
$ws.Range("J27").Value = 757.4167
$ws.Range("M27").Value = -793
$ws.Range("L110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("H110").Value = 0
$ws.Range("H126").Value = 3445.95
$ws.Range("I126").Value = 2455.2
$ws.Range("M126").Value = -4895.599999999999
$ws.Range("K126").Value = 7365.599999999999
$ws.Range("K132").Value = 4352.2104
$ws.Range("H132").Value = 1450.7368
$ws.Range("M132").Value = -1822.2104
$ws.Range("I132").Value = 1450.7368
$ws.Range("H136").Value = 3623.9473
$ws.Range("M136").Value = -9526.071599999999
$ws.Range("K136").Value = 12076.0716
$ws.Range("I136").Value = 4025.3572
$ws.Range("N110").ClearContents()

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2461.6667
$ws.Range("L113").Value = 15000
$ws.Range("I113").Value = 1192.5
$ws.Range("K113").Value = 3577.5
$ws.Range("J113").Value = 5000
$ws.Range("N113").Value = -19340
$ws.Range("M113").Value = -1407.5
$ws.Range("K132").Value = 3866.3748
$ws.Range("H132").Value = 1747.6666
$ws.Range("M132").Value = -1336.3748
$ws.Range("N132").Value = -15809.5001
$ws.Range("L132").Value = 10749.5001
$ws.Range("I132").Value = 1288.7916
$ws.Range("J132").Value = 3583.1667
